$d = $word.ActiveDocument
$rng = $d.Content
$rng.Collapse(0)
$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
<w:p><w:pPr><w:pStyle w:val="N3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:left="397" w:hanging="397"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="N3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:left="397" w:hanging="397"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="N3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:left="397" w:hanging="397"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="N3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:left="397" w:hanging="397"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="N3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:left="397" w:hanging="397"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:tab/></w:r><w:r><w:rPr><w:rStyle w:val="SigSignee"/></w:rPr><w:t>SIOBHIAN BROWN</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:rPr><w:rStyle w:val="Sigtitle"/></w:rPr><w:tab/><w:t>Authorised to sign by the Scottish Ministers</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SigBlock"/><w:rPr><w:rStyle w:val="Sigtitle"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="SigAdd"/></w:rPr><w:t>St Andrew’s House,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SigAdd"/></w:rPr><w:t>Edinburgh</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:rPr><w:rStyle w:val="SigDate"/></w:rPr><w:t>5th September 2024</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="N3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:left="397" w:hanging="397"/></w:pPr></w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
$rng.InsertXML($xml)
Write-Host "Inserted signature block paragraphs"
